# issue #5: add legislator_id, name, date into dataframe
# Adds three new columns (date, legislator_name, legislator_id) to the
# "股票" (stocks) worksheet, filling header + all data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$legislatorName = "鄭汝芬"
$legislatorId = 1713
$reportDate = "2011-11-21"

# --- header row (row 1) : copy the bold/bordered header style from G1 ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "date"

$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "legislator_name"

$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "legislator_id"

# --- data rows (2-11) ---
# The date string looks like a date, and Excel's normal cell-input parser
# (the same one backing plain ".Value =" assignment) would silently turn it
# into a date serial number. To keep it as literal text (matching the
# source data export), stage it via a formula (which yields a string
# result verbatim), then paste that *value* into the target cell, and
# finally paste the *format* from an existing plain data cell on top.
for ($r = 2; $r -le 11; $r++) {
    $helper = $ws.Range("Z" + (100 + $r))
    $helper.Formula = '="' + $reportDate + '"'
    $helper.Copy()
    $ws.Range("H" + $r).PasteSpecial(-4163)
    $ws.Range("C" + $r).Copy()
    $ws.Range("H" + $r).PasteSpecial(-4122)
    $helper.ClearContents()

    $ws.Range("C" + $r).Copy()
    $ws.Range("I" + $r).PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = $legislatorName

    $ws.Range("C" + $r).Copy()
    $ws.Range("J" + $r).PasteSpecial(-4122)
    $ws.Range("J" + $r).Value = $legislatorId
}
